# Update parameter files with WEP scaling and WFF_or_Ben
#
# Adds two new parameter rows to the "Parameters" sheet:
#   Row 56: MFTC_WEP_scaling = 1
#   Row 57: WFF_or_Benefit   = Max

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")

# --- Row 56: MFTC_WEP_scaling -------------------------------------------------
$ws.Range("B56").Value = "MFTC_WEP_scaling"

$ws.Range("C56").NumberFormat = "@"
$ws.Range("C56").Value = "1"

$ws.Range("D56").Value = "How should the Winter Energy Payment be scaled? Average week = 1, Winter week = 12/5, Summer week = 0"

# --- Row 57: WFF_or_Benefit ---------------------------------------------------
$ws.Range("B57").Value = "WFF_or_Benefit"

$ws.Range("C57").NumberFormat = "@"
$ws.Range("C57").Value = "Max"

$ws.Range("D57").Value = 'What work decision should we assume? Go off-benefit and receive IWTC = "WFF", stay on-benefit = "Benefit", or whichever gives a higher net income = "Max"'

# --- Formatting: match the style used for the rest of the parameter table ----
$newRows = $ws.Range("B56:D57")
$newRows.Interior.Pattern = 1
$newRows.Interior.Color = 15060409
$newRows.Font.Name = "Calibri"
$newRows.Font.Size = 11
$newRows.Font.Color = 0
$newRows.HorizontalAlignment = -4131
